$d = $word.ActiveDocument

# 1) editUser output: "user_id, username, full_name, status, and message"
#    -> "user_id, full_name, status, and message"
$d.Content.Find.Execute(", username, ", $true, $false, $false, $false, $false, $true, 1, $false, ", ", 2)

# 2) editContact output & getOneContact output both contain:
#    "contact_id, user_id, name, phone, address, website, email, status, and message"
#    -> "contact_id, name, phone, address, website, email, status, and message"
$d.Content.Find.Execute("contact_id, user_id,", $true, $false, $false, $false, $false, $true, 1, $false, "contact_id,", 2)

# 3) getOneContact required input: "contact_id and user_id" -> "contact_id"
$d.Content.Find.Execute("contact_id and user_id", $true, $false, $false, $false, $false, $true, 1, $false, "contact_id", 2)
